$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3012.322698989007
$ws.Range("C2").Value = 38.93031258922569
$ws.Range("D2").Value = 2355.426980468995
